# Generate Report for Handback
# Update the "handed back" timestamps once the localized files have come
# back in sync, across the Overview sheet and each per-locale sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
# "Latest HO Xliff Generate Date" for the second tracked file
# (4a28ded3-b3a0-46a7-a25c-d047914385c2.md) gets a fresh timestamp.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-03 00:52:22"

# --- zh-cn sheet -------------------------------------------------------
# Row 3 (4a28ded3-b3a0-46a7-a25c-d047914385c2 file) gets new handoff /
# handback datetimes now that the report has been regenerated.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-03 00:52:17"
$wsZhCn.Range("K3").Value = "2016-09-03 00:52:33"

# --- de-de sheet -------------------------------------------------------
# Row 3 (4a28ded3-b3a0-46a7-a25c-d047914385c2 file) likewise.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-03 00:52:22"
$wsDeDe.Range("K3").Value = "2016-09-03 00:52:40"
